$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.76548733449558
$ws.Range("D2").Value = 10.19202271113368
$ws.Range("E2").Value = 13.99772134997894
$ws.Range("F2").Value = 28.00002925415852
$ws.Range("G2").Value = 27.9926249224201
$ws.Range("H2").Value = 13.21923212800053
$ws.Range("J2").Value = 9.644818864116688
$ws.Range("N2").Value = 18.20530214750648
$ws.Range("O2").Value = 20.1862790208076
$ws.Range("B3").Value = 7.690488708080919
$ws.Range("D3").Value = 10.1879277601034
$ws.Range("E3").Value = 13.96069026936255
$ws.Range("F3").Value = 27.72584099711391
$ws.Range("G3").Value = 27.34581590149833
$ws.Range("H3").Value = 13.19068893700772
$ws.Range("J3").Value = 9.646505256136761
$ws.Range("N3").Value = 17.93039231733595
$ws.Range("O3").Value = 20.02620259791259
$ws.Range("B4").Value = 7.645832371061269
$ws.Range("D4").Value = 10.18718855004983
$ws.Range("E4").Value = 13.94095794876877
$ws.Range("F4").Value = 27.56376897813316
$ws.Range("G4").Value = 26.9496947766748
$ws.Range("H4").Value = 13.17586455310832
$ws.Range("J4").Value = 9.649197674905507
$ws.Range("N4").Value = 17.76180298315366
$ws.Range("O4").Value = 19.93251056327971
$ws.Range("B5").Value = 7.62800595514101
$ws.Range("D5").Value = 10.18733472167911
$ws.Range("E5").Value = 13.93367796296435
$ws.Range("F5").Value = 27.49937375959707
$ws.Range("G5").Value = 26.78879404821963
$ws.Range("H5").Value = 13.17050689977642
$ws.Range("J5").Value = 9.650711919510597
$ws.Range("N5").Value = 17.69323820309472
$ws.Range("O5").Value = 19.89552268896703
$ws.Range("B6").Value = 7.625068927627218
$ws.Range("D6").Value = 10.18738604175381
$ws.Range("E6").Value = 13.93251524037363
$ws.Range("F6").Value = 27.488782674454
$ws.Range("G6").Value = 26.76211597885145
$ws.Range("H6").Value = 13.16965865087502
$ws.Range("J6").Value = 9.65098855693036
$ws.Range("N6").Value = 17.68186382811895
$ws.Range("O6").Value = 19.8894539280562
$ws.Range("B7").Value = 7.645590426598107
$ws.Range("D7").Value = 10.18718870862609
$ws.Range("E7").Value = 13.94085668000986
$ws.Range("F7").Value = 27.56289374781946
$ws.Range("G7").Value = 26.94752234594259
$ws.Range("H7").Value = 13.17578952584707
$ws.Range("J7").Value = 9.649216407453391
$ws.Range("N7").Value = 17.76087762817026
$ws.Range("O7").Value = 19.93200685707592
$ws.Range("B8").Value = 7.739351513268095
$ws.Range("D8").Value = 10.19024303432427
$ws.Range("E8").Value = 13.98433288131076
$ws.Range("F8").Value = 27.90422929772455
$ws.Range("G8").Value = 27.76954980634348
$ws.Range("H8").Value = 13.20883245916042
$ws.Range("J8").Value = 9.645056714851505
$ws.Range("N8").Value = 18.11052099476146
$ws.Range("O8").Value = 20.13015363980903
$ws.Range("B9").Value = 7.9333151054759
$ws.Range("D9").Value = 10.2102635689738
$ws.Range("E9").Value = 14.09315794686801
$ws.Range("F9").Value = 28.6201392906549
$ws.Range("G9").Value = 29.37840297970911
$ws.Range("H9").Value = 13.29485483016452
$ws.Range("J9").Value = 9.650023774114835
$ws.Range("N9").Value = 18.79413998554022
$ws.Range("O9").Value = 20.55345977065718
$ws.Range("B10").Value = 8.080666255167973
$ws.Range("D10").Value = 10.23344268846874
$ws.Range("E10").Value = 14.1870649376271
$ws.Range("F10").Value = 29.16991955979584
$ws.Range("G10").Value = 30.54415329235382
$ws.Range("H10").Value = 13.37067982789178
$ws.Range("J10").Value = 9.661635036204377
$ws.Range("N10").Value = 19.2903645072747
$ws.Range("O10").Value = 20.88320280852609
$ws.Range("B11").Value = 8.14846916618764
$ws.Range("D11").Value = 10.24580427003773
$ws.Range("E11").Value = 14.23271206016586
$ws.Range("F11").Value = 29.42418676710619
$ws.Range("G11").Value = 31.06815545768676
$ws.Range("H11").Value = 13.40783664559342
$ws.Range("J11").Value = 9.66863450946393
$ws.Range("N11").Value = 19.51383106981896
$ws.Range("O11").Value = 21.036730358788
$ws.Range("B12").Value = 8.17423166240809
$ws.Range("D12").Value = 10.25074430124067
$ws.Range("E12").Value = 14.25040869153869
$ws.Range("F12").Value = 29.52098477406991
$ws.Range("G12").Value = 31.26546135151553
$ws.Range("H12").Value = 13.42228238153746
$ws.Range("J12").Value = 9.671530688125269
$ws.Range("N12").Value = 19.59804885977592
$ws.Range("O12").Value = 21.09532649437566
$ws.Range("B13").Value = 8.168679758780518
$ws.Range("D13").Value = 10.2496688993483
$ws.Range("E13").Value = 14.24657929015571
$ws.Range("F13").Value = 29.50011618201164
$ws.Range("G13").Value = 31.22302119215394
$ws.Range("H13").Value = 13.4191546736729
$ws.Range("J13").Value = 9.670896042929174
$ws.Range("N13").Value = 19.57993023154881
$ws.Range("O13").Value = 21.08268715005322
$ws.Range("B14").Value = 8.150587055476443
$ws.Range("D14").Value = 10.24620551355724
$ws.Range("E14").Value = 14.23415979427644
$ws.Range("F14").Value = 29.43214062403234
$ws.Range("G14").Value = 31.08441160933745
$ws.Range("H14").Value = 13.40901763072708
$ws.Range("J14").Value = 9.668867865316763
$ws.Range("N14").Value = 19.52076825227033
$ws.Range("O14").Value = 21.041542187005
$ws.Range("B15").Value = 8.139515373066374
$ws.Range("D15").Value = 10.24411774134096
$ws.Range("E15").Value = 14.22660571722196
$ws.Range("F15").Value = 29.3905677953432
$ws.Range("G15").Value = 30.99935692190976
$ws.Range("H15").Value = 13.40285702873039
$ws.Range("J15").Value = 9.667657493565205
$ws.Range("N15").Value = 19.48447494387068
$ws.Range("O15").Value = 21.01639798227277
$ws.Range("B16").Value = 8.076248655077
$ws.Range("D16").Value = 10.23267120236419
$ws.Range("E16").Value = 14.18413987375155
$ws.Range("F16").Value = 29.15337906325784
$ws.Range("G16").Value = 30.50976287739119
$ws.Range("H16").Value = 13.36830448962632
$ws.Range("J16").Value = 9.66121205336599
$ws.Range("N16").Value = 19.27570805562915
$ws.Range("O16").Value = 20.87323624657447
$ws.Range("B17").Value = 8.037617191298455
$ws.Range("D17").Value = 10.22611296865207
$ws.Range("E17").Value = 14.15883154443324
$ws.Range("F17").Value = 29.00887749313428
$ws.Range("G17").Value = 30.20763995635961
$ws.Range("H17").Value = 13.34778449085651
$ws.Range("J17").Value = 9.657696969332493
$ws.Range("N17").Value = 19.14699881874138
$ws.Range("O17").Value = 20.78627997687077
$ws.Range("B18").Value = 8.015471503473869
$ws.Range("D18").Value = 10.22251205291303
$ws.Range("E18").Value = 14.14455113750956
$ws.Range("F18").Value = 28.92616169188517
$ws.Range("G18").Value = 30.03328550724615
$ws.Range("H18").Value = 13.33623314540072
$ws.Range("J18").Value = 9.655836931840597
$ws.Range("N18").Value = 19.07275961144547
$ws.Range("O18").Value = 20.73659984591649
$ws.Range("B19").Value = 8.007986816985902
$ws.Range("D19").Value = 10.22132231765953
$ws.Range("E19").Value = 14.13976378266686
$ws.Range("F19").Value = 28.89822649347396
$ws.Range("G19").Value = 29.97415906526007
$ws.Range("H19").Value = 13.33236543278866
$ws.Range("J19").Value = 9.65523497299491
$ws.Range("N19").Value = 19.04759001474562
$ws.Range("O19").Value = 20.71983797588696
$ws.Range("B20").Value = 8.041722081226208
$ws.Range("D20").Value = 10.22679340089651
$ws.Range("E20").Value = 14.16149713654979
$ws.Range("F20").Value = 29.02421940351396
$ws.Range("G20").Value = 30.23986323467787
$ws.Range("H20").Value = 13.34994293450565
$ws.Range("J20").Value = 9.658054425068778
$ws.Range("N20").Value = 19.16072234871255
$ws.Range("O20").Value = 20.79550231663107
$ws.Range("B21").Value = 8.155899150234889
$ws.Range("D21").Value = 10.24721578547873
$ws.Range("E21").Value = 14.23779662850099
$ws.Range("F21").Value = 29.45209347683156
$ws.Range("G21").Value = 31.12515671041965
$ws.Range("H21").Value = 13.41198500623132
$ws.Range("J21").Value = 9.669456935443925
$ws.Range("N21").Value = 19.53815713259453
$ws.Range("O21").Value = 21.05361540318226
$ws.Range("B22").Value = 8.231016616546755
$ws.Range("D22").Value = 10.26207120930414
$ws.Range("E22").Value = 14.29005360826733
$ws.Range("F22").Value = 29.7346810838264
$ws.Range("G22").Value = 31.69711573386239
$ws.Range("H22").Value = 13.45471638674811
$ws.Range("J22").Value = 9.678340111133934
$ws.Range("N22").Value = 19.78244317294327
$ws.Range("O22").Value = 21.22495736481477
$ws.Range("B23").Value = 8.190887495042949
$ws.Range("D23").Value = 10.25400542471188
$ws.Range("E23").Value = 14.26194779235945
$ws.Range("F23").Value = 29.58361782623343
$ws.Range("G23").Value = 31.39252465664353
$ws.Range("H23").Value = 13.43171277717938
$ws.Range("J23").Value = 9.67346854575416
$ws.Range("N23").Value = 19.65230661567169
$ws.Range("O23").Value = 21.13328259258936
$ws.Range("B24").Value = 8.039866059139614
$ws.Range("D24").Value = 10.22648524933156
$ws.Range("E24").Value = 14.16029118155904
$ws.Range("F24").Value = 29.01728219750546
$ws.Range("G24").Value = 30.22529712948333
$ws.Range("H24").Value = 13.34896633560851
$ws.Range("J24").Value = 9.65789231825727
$ws.Range("N24").Value = 19.15451869034849
$ws.Range("O24").Value = 20.79133192062987
$ws.Range("B25").Value = 7.879895846024594
$ws.Range("D25").Value = 10.2033530996724
$ws.Range("E25").Value = 14.0612345977308
$ws.Range("F25").Value = 28.42194760795594
$ws.Range("G25").Value = 28.94502986164271
$ws.Range("H25").Value = 13.26934136391781
$ws.Range("J25").Value = 9.647278360295342
$ws.Range("N25").Value = 18.60992361498239
$ws.Range("O25").Value = 20.435470315409
